# Apply the "cryptos list" data refresh described by the commit:
#   "Updated cryptos list on Wed Apr 12 08:57:20 UTC 2023 with GitHub Actions"
#
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Most rows only get new
# Price/Volume figures, but a handful of rows also shift which coin occupies
# them (rows 40-42 rotate, and a new row 51 entry replaces the old one).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns are plain text in the workbook (e.g. "91.80",
# "1.001", percentages with padding spaces). Force the range to Text format
# first so Excel does not "helpfully" reinterpret numeric-looking strings
# as numbers (which would drop trailing zeros / change precision).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$updates = @{
    'D2' = '30.006.48'
    'E2' = '  -0.41%  '
    'D3' = '1.872.16'
    'E3' = '  -2.53%  '
    'E4' = '  +0.09%  '
    'D5' = '319.57'
    'E5' = '  -3.29%  '
    'E6' = '  +0.05%  '
    'D7' = '0.5040'
    'E7' = '  -3.18%  '
    'D8' = '0.3947'
    'D9' = '0.08198'
    'E9' = '  -3.81%  '
    'D10' = '42.19'
    'E10' = '  -2.53%  '
    'E11' = '  -3.23%  '
    'D12' = '23.54'
    'E12' = '  +5.26%  '
    'D13' = '1.869.37'
    'E13' = '  -2.78%  '
    'D14' = '6.289'
    'E14' = '  -1.97%  '
    'D15' = '7.182'
    'E15' = '  -3.10%  '
    'E16' = '  +0.06%  '
    'D17' = '91.80'
    'E17' = '  -4.08%  '
    'D18' = '0.00001086'
    'E18' = '  -2.51%  '
    'D19' = '0.06407'
    'E19' = '  -4.68%  '
    'D20' = '18.12'
    'E20' = '  -1.05%  '
    'E21' = '  +0.00%  '
    'D22' = '30.023.34'
    'E22' = '  -0.36%  '
    'D23' = '5.843'
    'E23' = '  -3.12%  '
    'E24' = '  -2.03%  '
    'D25' = '2.172'
    'E25' = '  -2.30%  '
    'D26' = '2.082.47'
    'E26' = '  -2.84%  '
    'D27' = '21.29'
    'E27' = '  +0.70%  '
    'D28' = '160.53'
    'E28' = '  +0.36%  '
    'D29' = '2.217'
    'E29' = '  -9.76%  '
    'D30' = '127.18'
    'E30' = '  -1.44%  '
    'D31' = '1.067'
    'E31' = '  -0.97%  '
    'D32' = '0.1033'
    'E32' = '  -2.12%  '
    'E33' = '  -2.57%  '
    'D34' = '3.679'
    'E34' = '  +1.27%  '
    'D35' = '0.02436'
    'E35' = '  -2.46%  '
    'D36' = '5.207'
    'E36' = '  -0.19%  '
    'D37' = '0.06357'
    'E37' = '  -3.87%  '
    'D38' = '0.2141'
    'E38' = '  -3.47%  '
    'D39' = '1.172'
    'E39' = '  -5.18%  '
    'B40' = 'TrustWalletToken'
    'C40' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D40' = '1.221'
    'E40' = '  -2.00%  '
    'B41' = 'FraxShare'
    'C41' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D41' = '8.474'
    'E41' = '  -5.22%  '
    'B42' = 'TheSandbox'
    'C42' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D42' = '0.6308'
    'E42' = '  -3.47%  '
    'D43' = '11.27'
    'E43' = '  -3.51%  '
    'D44' = '1.001'
    'E44' = '  +0.01%  '
    'D45' = '0.5904'
    'D46' = '12.90'
    'E46' = '  -3.01%  '
    'D47' = '2.094'
    'E47' = '  +0.21%  '
    'E48' = '  -3.88%  '
    'D49' = '122.65'
    'E49' = '  -1.73%  '
    'D50' = '1.203'
    'E50' = '  -3.84%  '
    'B51' = 'WEMIXTOKEN'
    'C51' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D51' = '1.123'
    'E51' = '  -3.45%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Drop the temporary Text number format again so the touched cells end up
# with the same (default/general) style they started with.
$dataRange.ClearFormats()
